# Add ", APIs" right after the existing "Heroku" text in the Technical
# Skills line, matching the run formatting already used on that line, and
# move the "_GoBack" bookmark from the end of the "...ever invented."
# paragraph to the end of this line (after the newly inserted text).

$d = $word.ActiveDocument

# --- 1) Insert ", APIs" immediately after "Heroku" -------------------------
$rng = $d.Content
$rng.Find.Execute("Heroku", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Collapse(0)

# Append a temporary trailing marker character along with the real text so
# the insertion point we need for the bookmark is never the very last
# character position of the paragraph (collapsing a range exactly at a
# paragraph's end confuses bookmark placement) -- we strip the marker off
# again below once the bookmark has been anchored safely.
$rng.InsertAfter(", APIsZ")

# Match the surrounding run formatting (rFonts ascii/eastAsia/hAnsi/cs +
# color=auto), same as the rest of the "Technical Skills" line.
$rng.Font.NameAscii = "Cambria"
$rng.Font.NameFarEast = "Georgia"
$rng.Font.NameOther = "Cambria"
$rng.Font.NameBi = "Calibri"
$rng.Font.Color = -16777216

# --- 2) Move the "_GoBack" bookmark to sit right after ", APIs" ------------
# Wrap the trailing marker character with a (non-collapsed) range, anchor
# the bookmark there -- re-adding "_GoBack" automatically removes it from
# its previous location -- then delete the marker so the bookmark collapses
# to an empty range right after "APIs", exactly like the original markup.
$marker = $d.Range($rng.End - 1, $rng.End)
$d.Bookmarks.Add("_GoBack", $marker)
$bm = $d.Bookmarks.Item("_GoBack")
$bm.Range.Text = ""

Write-Output "done"
